$d = $word.ActiveDocument

# --- Table 2 (TABELA 2 - MATERIAL ENCAMINHADO A EXAME), data row ---
$t2 = $d.Tables.Item(2)
$t2.Cell(3, 4).Range.Text = "BLAZER"   # Dito no oficio: A USA -> BLAZER
$t2.Cell(3, 5).Range.Text = "654"      # Lacre de Entrada: 345 -> 654

# --- Section heading "3.2 DOS CARTUCHOS " -> "3.1 DOS CARTUCHOS " ---
$d.Content.Find.Execute("3.2 DOS CARTUCHOS ", $true, $false, $false, $false, $false, $true, 1, $false, "3.1 DOS CARTUCHOS ", 2)

# --- Table 4 (TABELA 4 - DESCRICAO DO(S) CARTUCHO(S) INTACTOS), data row ---
$t4 = $d.Tables.Item(4)
$t4.Cell(3, 3).Range.Text = ".380 AUTO"        # Calibre Nominal: .22 Curto -> .380 AUTO
$t4.Cell(3, 4).Range.Text = "BLAZER"           # Marca: A USA -> BLAZER
$t4.Cell(3, 5).Range.Text = "NORTE AMERICANA"  # Procedencia: ESTADUNIDENSE -> NORTE AMERICANA
$t4.Cell(3, 6).Range.Text = "AÇO"              # Espoleta: NIQUELADA -> ACO
$t4.Cell(3, 7).Range.Text = "ALUMÍNIO"         # Estojo (Lote): NIQUELADO -> ALUMINIO
$t4.Cell(3, 8).Range.Text = "CSCV"             # Projetil: CHPP -> CSCV

# --- Legend paragraph inside the "Legenda:" cell (row 4, col 1) ---
$legendCell = $t4.Cell(4, 1)
$legendCell.Range.Find.Execute("CHPP  Chumbo Ponta Plana", $true, $false, $false, $false, $false, $true, 1, $false, "CSCV  Chumbo Semi Canto Vivo ", 2)
